$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Apparel and Clothing, 13, 681
$ws.Range("A7").Value = "Apparel and Clothing"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 681

# Update the selected cell/range in the sheet view
$ws.Range("H11").Select() | Out-Null
